$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("main")
$ws.Range("A38").Value = "test"
